# New baseline: new FGS and AIRS loaders, polynomial drift
#
# Applies to the "Submissions" sheet (internal sheet2.xml):
#  - Apply a "0.000" number format to the whole "LB" (column N) data range
#  - Add five new submission rows (36-40): new FGS / AIRS loaders + poly drift
#  - Grow the Table2 listobject + sheet dimension/selection to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submissions")
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Apply the "0.000" number format to the existing LB column (N2:N35),
#    including the rows that don't yet hold a value - only the format is
#    touched here, the stored values themselves are left alone.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 14).NumberFormat = "0.000"
}

# ---------------------------------------------------------------------
# 2. New rows appended to the submissions log.
# ---------------------------------------------------------------------
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 4).Value = "(29)+no sanity checks"
$ws.Cells.Item(36, 14).Value = 0.599
$ws.Cells.Item(36, 14).NumberFormat = "0.000"

$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 4).Value = "(29)+new FGS"
$ws.Cells.Item(37, 2).Value = "0a3b"
$ws.Cells.Item(37, 3).Value = "Yes"
$ws.Cells.Item(37, 5).Value = "Yes"
$ws.Cells.Item(37, 7).Value = 0.6079
$ws.Cells.Item(37, 8).Value = 199.05
$ws.Cells.Item(37, 9).Value = 151.75
$ws.Cells.Item(37, 14).Value = 0.602
$ws.Cells.Item(37, 14).NumberFormat = "0.000"

$ws.Cells.Item(38, 4).Value = "(35)+fixed FGS mean"
$ws.Cells.Item(38, 5).Value = "Yes"
$ws.Cells.Item(38, 7).Value = 0.6077
$ws.Cells.Item(38, 8).Value = 196.7
$ws.Cells.Item(38, 9).Value = 148.9
$ws.Cells.Item(38, 14).NumberFormat = "0.000"

$ws.Cells.Item(39, 4).Value = "(36)+poly FGS"
$ws.Cells.Item(39, 5).Value = "Yes"
$ws.Cells.Item(39, 7).Value = 0.6071
$ws.Cells.Item(39, 8).Value = 228.8
$ws.Cells.Item(39, 9).Value = 173.5
$ws.Cells.Item(39, 14).Value = 0.6
$ws.Cells.Item(39, 14).NumberFormat = "0.000"

$ws.Cells.Item(40, 4).Value = "(37)+new AIRS simpler noise_est"
$ws.Cells.Item(40, 5).Value = "Yes"
$ws.Cells.Item(40, 7).Value = 0.6072
$ws.Cells.Item(40, 8).Value = 232.3
$ws.Cells.Item(40, 9).Value = 178.4
$ws.Cells.Item(40, 14).Value = 0.599
$ws.Cells.Item(40, 14).NumberFormat = "0.000"

# ---------------------------------------------------------------------
# 3. Grow Table2 to cover the new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("A1:O40"))

# ---------------------------------------------------------------------
# 4. Scroll/selection bookkeeping to match where the author left off.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("N41").Select()

Write-Output "applied submissions baseline update"
